# Regenerate merged AHB files
#
# 1) Rename the header row (row 1) from the "_old"/"_new" variant suffixes
#    to the "_FV2404"/"_FV2410" release-tag suffixes (columns A:J and L:U;
#    column K, "diff", is left untouched).
# 2) Turn the A1:U84 range into a proper Excel Table ("Table1") using the
#    same (renamed) headers, row-stripe styling and an autofilter.
# 3) Freeze the header row (split/freeze below row 1) and leave the
#    selection on A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2404", "Segmentgruppe_FV2404", "Segment_FV2404", "Datenelement_FV2404", "Segment ID_FV2404",
    "Code_FV2404", "Qualifier_FV2404", "Beschreibung_FV2404", "Bedingungsausdruck_FV2404", "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410", "Segmentgruppe_FV2410", "Segment_FV2410", "Datenelement_FV2410", "Segment ID_FV2410",
    "Code_FV2410", "Qualifier_FV2410", "Beschreibung_FV2410", "Bedingungsausdruck_FV2410", "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Convert the used range into an Excel Table ("ListObject"), matching the
# workbook's ref="A1:U84" table definition.
$rng = $ws.Range("A1:U84")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# Freeze panes below row 1 (split at row 2), then restore the selection to A1.
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
